$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 463652.7
$ws.Range("J17").Value = 514995.7
$ws.Range("L17").Value = 1544987.1
$ws.Range("N17").Value = -1545323.1
$ws.Range("H40").Value = 1151.6666
$ws.Range("I40").Value = 1110.4736
$ws.Range("K40").Value = 1110.4736
$ws.Range("M40").Value = -935.4736
$ws.Range("H76").Value = 3527.5715
$ws.Range("I76").Value = 3623.5
$ws.Range("J76").Value = 3399.6667
$ws.Range("K76").Value = 3623.5
$ws.Range("L76").Value = 3399.6667
$ws.Range("M76").Value = -3308.5
$ws.Range("N76").Value = -4029.6667
$ws.Range("H79").Value = 3527.5715
$ws.Range("I79").Value = 3623.5
$ws.Range("J79").Value = 3399.6667
$ws.Range("K79").Value = 3623.5
$ws.Range("L79").Value = 3399.6667
$ws.Range("M79").Value = -2531.5
$ws.Range("N79").Value = -5583.6667
$ws.Range("H86").Value = 3469.8235
$ws.Range("I86").Value = 2806.7693
$ws.Range("J86").Value = 5624.75
$ws.Range("K86").Value = 2806.7693
$ws.Range("L86").Value = 5624.75
$ws.Range("M86").Value = -1683.7693
$ws.Range("N86").Value = -7870.75
$ws.Range("H89").Value = 3469.8235
$ws.Range("I89").Value = 2806.7693
$ws.Range("J89").Value = 5624.75
$ws.Range("K89").Value = 14033.8465
$ws.Range("L89").Value = 28123.75
$ws.Range("M89").Value = -8417.8465
$ws.Range("N89").Value = -39355.75
$ws.Range("H99").Value = 6575
$ws.Range("I99").Value = 400
$ws.Range("J99").Value = 12750
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 38250
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -41246
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H132").Value = 56679.5
$ws.Range("I132").Value = 30569.47
$ws.Range("K132").Value = 91708.41
$ws.Range("M132").Value = -89178.41
$ws.Range("H133").Value = 68500
$ws.Range("J133").Value = 68500
$ws.Range("L133").Value = 68500
$ws.Range("N133").Value = -78620
$ws.Range("H137").Value = 1669
$ws.Range("I137").Value = 1742
$ws.Range("J137").Value = 1450
$ws.Range("K137").Value = 5226
$ws.Range("L137").Value = 4350
$ws.Range("M137").Value = -2676
$ws.Range("N137").Value = -9450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2133.4644
$ws.Range("I32").Value = 1778.091
$ws.Range("J32").Value = 6042.5713
$ws.Range("K32").Value = 1778.091
$ws.Range("L32").Value = 6042.5713
$ws.Range("M32").Value = -1491.091
$ws.Range("N32").Value = -6616.5713
$ws.Range("H74").Value = 3888.6453
$ws.Range("I74").Value = 483.37036
$ws.Range("K74").Value = 483.37036
$ws.Range("M74").Value = 390.62964
$ws.Range("H77").Value = 3888.6453
$ws.Range("I77").Value = 483.37036
$ws.Range("K77").Value = 2416.8518
$ws.Range("M77").Value = 1951.1482
$ws.Range("H88").Value = 2707.2
$ws.Range("I88").Value = 2765
$ws.Range("J88").Value = 2682.4285
$ws.Range("K88").Value = 2765
$ws.Range("L88").Value = 2682.4285
$ws.Range("M88").Value = -2359
$ws.Range("N88").Value = -3494.4285
$ws.Range("H91").Value = 2707.2
$ws.Range("I91").Value = 2765
$ws.Range("J91").Value = 2682.4285
$ws.Range("K91").Value = 2765
$ws.Range("L91").Value = 2682.4285
$ws.Range("M91").Value = -1361
$ws.Range("N91").Value = -5490.4285
$ws.Range("H102").Value = 3691.3
$ws.Range("J102").Value = 3666.6667
$ws.Range("L102").Value = 3666.6667
$ws.Range("N102").Value = -6910.6667
$ws.Range("H106").Value = 196799.8
$ws.Range("J106").Value = 196799.8
$ws.Range("L106").Value = 196799.8
$ws.Range("N106").Value = -199323.8
$ws.Range("H132").Value = 1875.3226
$ws.Range("I132").Value = 1659.8276
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4979.4828
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2449.4828
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H101").Value = 20000
$ws.Range("I101").Value = 20000
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 20000
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -16755
$ws.Range("N101").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = 0
$ws.Range("H107").Value = 1685.9375
$ws.Range("I107").Value = 1575.3636
$ws.Range("J107").Value = 1929.2
$ws.Range("K107").Value = 1575.3636
$ws.Range("L107").Value = 1929.2
$ws.Range("M107").Value = 344.6364000000001
$ws.Range("N107").Value = -5769.2
$ws.Range("H134").Value = 2834.3
$ws.Range("I134").Value = 2452.6
$ws.Range("K134").Value = 7357.799999999999
$ws.Range("M134").Value = -4822.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30062.514
$ws.Range("I31").Value = 32750.562
$ws.Range("K31").Value = 32750.562
$ws.Range("M31").Value = -32455.562
$ws.Range("H34").Value = 30062.514
$ws.Range("I34").Value = 32750.562
$ws.Range("K34").Value = 32750.562
$ws.Range("M34").Value = -32548.562
$ws.Range("H58").Value = 2924
$ws.Range("I58").Value = 3066.1667
$ws.Range("K58").Value = 3066.1667
$ws.Range("M58").Value = -2863.1667
$ws.Range("H97").Value = 40001
$ws.Range("J97").Value = 40001
$ws.Range("L97").Value = 40001
$ws.Range("N97").Value = -41983
$ws.Range("H107").Value = 834.6087
$ws.Range("I107").Value = 950.86664
$ws.Range("J107").Value = 616.625
$ws.Range("K107").Value = 950.86664
$ws.Range("L107").Value = 616.625
$ws.Range("M107").Value = 969.13336
$ws.Range("N107").Value = -4456.625
$ws.Range("H134").Value = 17808.234
$ws.Range("I134").Value = 8169.6665
$ws.Range("J134").Value = 54985.57
$ws.Range("K134").Value = 24508.9995
$ws.Range("L134").Value = 164956.71
$ws.Range("M134").Value = -21973.9995
$ws.Range("N134").Value = -170026.71
$ws.Range("H136").Value = 2924
$ws.Range("I136").Value = 3066.1667
$ws.Range("K136").Value = 9198.500100000001
$ws.Range("M136").Value = -6648.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9232.5
$ws.Range("I5").Value = 600
$ws.Range("J5").Value = 13548.75
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 40646.25
$ws.Range("M5").Value = -1688
$ws.Range("N5").Value = -40870.25
$ws.Range("H14").Value = 286597
$ws.Range("I14").Value = 286597
$ws.Range("K14").Value = 859791
$ws.Range("M14").Value = -859618
$ws.Range("H107").Value = 790.94446
$ws.Range("I107").Value = 659.5
$ws.Range("J107").Value = 1251
$ws.Range("K107").Value = 1978.5
$ws.Range("L107").Value = 3753
$ws.Range("M107").Value = -58.5
$ws.Range("N107").Value = -7593
$ws.Range("H122").Value = 2694.3572
$ws.Range("I122").Value = 2471.5
$ws.Range("J122").Value = 2991.5
$ws.Range("K122").Value = 22243.5
$ws.Range("L122").Value = 26923.5
$ws.Range("M122").Value = -19793.5
$ws.Range("N122").Value = -31823.5
$ws.Range("H135").Value = 9232.5
$ws.Range("I135").Value = 600
$ws.Range("J135").Value = 13548.75
$ws.Range("K135").Value = 5400
$ws.Range("L135").Value = 121938.75
$ws.Range("M135").Value = -2865
$ws.Range("N135").Value = -127008.75
$ws.Range("H136").Value = 1122532
$ws.Range("I136").Value = 3334273
$ws.Range("J136").Value = 16661.5
$ws.Range("K136").Value = 10002819
$ws.Range("L136").Value = 49984.5
$ws.Range("M136").Value = -9997719
$ws.Range("N136").Value = -60184.5
$ws.Range("H137").Value = 4326
$ws.Range("I137").Value = 2790.7778
$ws.Range("K137").Value = 8372.3334
$ws.Range("M137").Value = -3272.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2402
$ws.Range("I80").Value = 2190.4285
$ws.Range("J80").Value = 2613.5715
$ws.Range("K80").Value = 2190.4285
$ws.Range("L80").Value = 2613.5715
$ws.Range("M80").Value = -1192.4285
$ws.Range("N80").Value = -4609.5715
$ws.Range("H83").Value = 2402
$ws.Range("I83").Value = 2190.4285
$ws.Range("J83").Value = 2613.5715
$ws.Range("K83").Value = 10952.1425
$ws.Range("L83").Value = 13067.8575
$ws.Range("M83").Value = -5960.1425
$ws.Range("N83").Value = -23051.8575
$ws.Range("H122").Value = 2394.8
$ws.Range("I122").Value = 2383.7273
$ws.Range("J122").Value = 2425.25
$ws.Range("K122").Value = 7151.1819
$ws.Range("L122").Value = 7275.75
$ws.Range("M122").Value = -4701.1819
$ws.Range("N122").Value = -12175.75
$ws.Range("H132").Value = 253389.4
$ws.Range("I132").Value = 253389.4
$ws.Range("K132").Value = 760168.2
$ws.Range("M132").Value = -757638.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 148.13637
$ws.Range("I55").Value = 147.88889
$ws.Range("J55").Value = 149.25
$ws.Range("K55").Value = 147.88889
$ws.Range("L55").Value = 149.25
$ws.Range("M55").Value = 25.11111
$ws.Range("N55").Value = -495.25
$ws.Range("H68").Value = 1275
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1275
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 1275
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -2773
$ws.Range("H71").Value = 1275
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1275
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 6375
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -13863
$ws.Range("H100").Value = 62151.945
$ws.Range("I100").Value = 67483.19
$ws.Range("J100").Value = 19502
$ws.Range("K100").Value = 67483.19
$ws.Range("L100").Value = 19502
$ws.Range("M100").Value = -66942.19
$ws.Range("N100").Value = -20584
$ws.Range("H133").Value = 51316.332
$ws.Range("J133").Value = 51316.332
$ws.Range("L133").Value = 51316.332
$ws.Range("N133").Value = -56376.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2284.4
$ws.Range("I107").Value = 1023.5
$ws.Range("J107").Value = 3125
$ws.Range("K107").Value = 3070.5
$ws.Range("L107").Value = 9375
$ws.Range("M107").Value = -1150.5
$ws.Range("N107").Value = -13215
$ws.Range("H132").Value = 3975.1853
$ws.Range("I132").Value = 3975.1853
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11925.5559
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9395.555899999999
$ws.Range("N132").ClearContents()
